# NCI Thesaurus, Mitelmandb update - aug 2025
$wb = $excel.ActiveWorkbook

# --- "compounds" sheet: update NCI Thesaurus source_version ---
$wsCompounds = $wb.Worksheets.Item("compounds")
$wsCompounds.Range("E3").Value = "25.07d"

# --- "biomarkers" sheet: update Mitelman Database source_version ---
$wsBiomarkers = $wb.Worksheets.Item("biomarkers")
$wsBiomarkers.Range("E3").Value = "v20250710"

# Update the remembered selection on the biomarkers sheet (E6 -> E3), then
# restore "compounds" as the active/selected sheet tab.
$wsBiomarkers.Range("E3").Select()
$wsCompounds.Activate()
